$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D22").Value = 0.6447426901493167
$ws.Range("C23").Value = 0.2386249091493167
$ws.Range("D23").Value = 0.597740902
$ws.Range("B24").Value = -0.0107480648506833
$ws.Range("C24").Value = 0.042359665
